$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.988.24"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "3.377.63"
$ws.Range("E3").Value = "  +7.90%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'259.35"
$ws.Range("E5").Value = "  +7.96%  "
$ws.Range("D6").Value = "'628.89"
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("E7").Value = "  +25.49%  "
$ws.Range("D8").Value = "'0.394"
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'0.878"
$ws.Range("E10").Value = "  +12.24%  "
$ws.Range("D11").Value = "3.374.35"
$ws.Range("E11").Value = "  +7.88%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "98.748.60"
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("D14").Value = "'36.11"
$ws.Range("E14").Value = "  +6.14%  "
$ws.Range("D15").Value = "'0.0000249"
$ws.Range("E15").Value = "  +3.08%  "
$ws.Range("D16").Value = "3.962.29"
$ws.Range("E16").Value = "  +6.84%  "
$ws.Range("E17").Value = "  +1.66%  "
$ws.Range("D18").Value = "3.363.73"
$ws.Range("E18").Value = "  +7.72%  "
$ws.Range("D19").Value = "'3.57"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").Value = "'15.27"
$ws.Range("E20").Value = "  +4.91%  "
$ws.Range("D21").Value = "'493.74"
$ws.Range("E21").Value = "  -6.65%  "
$ws.Range("D22").Value = "'6.18"
$ws.Range("E22").Value = "  +8.37%  "
$ws.Range("E23").Value = "  +9.42%  "
$ws.Range("E24").Value = "  +7.01%  "
$ws.Range("D25").Value = "'5.67"
$ws.Range("E25").Value = "  +3.54%  "
$ws.Range("D26").Value = "'88.77"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'11.98"
$ws.Range("E27").Value = "  +3.22%  "
$ws.Range("D29").Value = "'0.281"
$ws.Range("E29").Value = "  +18.79%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "'0.195"
$ws.Range("E31").Value = "  +11.62%  "
$ws.Range("D32").Value = "'0.133"
$ws.Range("E32").Value = "  +5.37%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +16.38%  "
$ws.Range("E34").Value = "  +6.37%  "
$ws.Range("D35").Value = "'27.84"
$ws.Range("E35").Value = "  +4.01%  "
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  +5.33%  "
$ws.Range("D39").Value = "'0.464"
$ws.Range("E39").Value = "  +5.37%  "
$ws.Range("D40").Value = "'499.36"
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("D41").Value = "'24.88"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("D42").Value = "'3.82"
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("E43").Value = "  +4.05%  "
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").Value = "'0.789"
$ws.Range("E45").Value = "  +12.41%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'160.56"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "'1.95"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").Value = "'0.835"
$ws.Range("E49").Value = "  +14.24%  "
$ws.Range("E50").Value = "  +4.01%  "
$ws.Range("D51").Value = "'46.12"
$ws.Range("E51").Value = "  +3.93%  "
